# Update gh-pages to output generated at 456a3b4
# This applies the refreshed "想去人数" (want-to-go count) figures captured
# in column F across the four worksheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 200
$ws1.Range("F4").Value = 675
$ws1.Range("F7").Value = 1299
$ws1.Range("F8").Value = 771
$ws1.Range("F11").Value = 2782
$ws1.Range("F17").Value = 852
$ws1.Range("F18").Value = 71
$ws1.Range("F21").Value = 90
$ws1.Range("F27").Value = 4840
$ws1.Range("F29").Value = 115

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F27").Value = 506
$ws2.Range("F37").Value = 675

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F6").Value = 357
$ws3.Range("F7").Value = 328

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 200
$ws4.Range("F6").Value = 357
$ws4.Range("F9").Value = 675
$ws4.Range("F14").Value = 1299
$ws4.Range("F15").Value = 771
$ws4.Range("F20").Value = 2782
$ws4.Range("F21").Value = 2782
$ws4.Range("F28").Value = 328
$ws4.Range("F30").Value = 852
$ws4.Range("F31").Value = 852
$ws4.Range("F32").Value = 71
$ws4.Range("F36").Value = 90
$ws4.Range("F42").Value = 506
$ws4.Range("F46").Value = 4840
$ws4.Range("F50").Value = 115
$ws4.Range("F51").Value = 675
